# Applies the "Automatic update of files." diff to the Artfynd worksheet.
# Every target cell's final value is known exactly from the diff, so we
# just assign the new literal value to each cell (using Value2, which
# round-trips numbers/strings correctly through this COM shim).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value2 = 79002

# Row 3
$ws.Cells.Item(3, 1).Value2  = 131033355
$ws.Cells.Item(3, 2).Value2  = 79244
$ws.Cells.Item(3, 17).Value2 = 396033   # Q3
$ws.Cells.Item(3, 18).Value2 = 6804799  # R3
$ws.Cells.Item(3, 26).Value2 = "14:50"  # Z3
$ws.Cells.Item(3, 28).Value2 = "14:50"  # AB3

# Row 4
$ws.Cells.Item(4, 1).Value2  = 131033337
$ws.Cells.Item(4, 2).Value2  = 79244
$ws.Cells.Item(4, 17).Value2 = 395826   # Q4
$ws.Cells.Item(4, 18).Value2 = 6804765  # R4
$ws.Cells.Item(4, 26).Value2 = "13:23"  # Z4
$ws.Cells.Item(4, 28).Value2 = "13:23"  # AB4

# Row 5
$ws.Cells.Item(5, 2).Value2 = 78647

# Row 6
$ws.Cells.Item(6, 2).Value2 = 79244

# Row 8
$ws.Cells.Item(8, 1).Value2  = 131033281
$ws.Cells.Item(8, 2).Value2  = 79863
$ws.Cells.Item(8, 5).Value2  = 6453                               # E8
$ws.Cells.Item(8, 6).Value2  = "Vedskivlav"                       # F8
$ws.Cells.Item(8, 7).Value2  = "Hertelidea botryosa"              # G8
$ws.Cells.Item(8, 8).Value2  = "(Fr.) Printzen & Kantvilas"       # H8
$ws.Cells.Item(8, 17).Value2 = 395662                             # Q8
$ws.Cells.Item(8, 18).Value2 = 6804783                            # R8
$ws.Cells.Item(8, 26).Value2 = "11:57"                            # Z8
$ws.Cells.Item(8, 28).Value2 = "11:57"                            # AB8

# Row 9
$ws.Cells.Item(9, 1).Value2  = 131033360
$ws.Cells.Item(9, 2).Value2  = 78910
$ws.Cells.Item(9, 5).Value2  = 353                                # E9
$ws.Cells.Item(9, 6).Value2  = "Dvärgbägarlav"                    # F9
$ws.Cells.Item(9, 7).Value2  = "Cladonia parasitica"              # G9
$ws.Cells.Item(9, 8).Value2  = "(Hoffm.) Hoffm."                  # H9
$ws.Cells.Item(9, 17).Value2 = 395791                             # Q9
$ws.Cells.Item(9, 18).Value2 = 6804722                            # R9
$ws.Cells.Item(9, 26).Value2 = "13:19"                            # Z9
$ws.Cells.Item(9, 28).Value2 = "13:19"                            # AB9

# Row 10
$ws.Cells.Item(10, 2).Value2 = 79244

# Row 11
$ws.Cells.Item(11, 2).Value2 = 78910

# Row 12
$ws.Cells.Item(12, 2).Value2 = 79244

# Row 13
$ws.Cells.Item(13, 1).Value2  = 131033319
$ws.Cells.Item(13, 2).Value2  = 79244
$ws.Cells.Item(13, 17).Value2 = 395654   # Q13
$ws.Cells.Item(13, 18).Value2 = 6804639  # R13
$ws.Cells.Item(13, 26).Value2 = "12:01"  # Z13
$ws.Cells.Item(13, 28).Value2 = "12:01"  # AB13

# Row 14
$ws.Cells.Item(14, 1).Value2  = 131033320
$ws.Cells.Item(14, 2).Value2  = 79244
$ws.Cells.Item(14, 17).Value2 = 395688   # Q14
$ws.Cells.Item(14, 18).Value2 = 6804619  # R14
$ws.Cells.Item(14, 26).Value2 = "12:03"  # Z14
$ws.Cells.Item(14, 28).Value2 = "12:03"  # AB14

# Row 17
$ws.Cells.Item(17, 2).Value2 = 79244

# Row 18
$ws.Cells.Item(18, 2).Value2 = 79244

# Row 19
$ws.Cells.Item(19, 2).Value2 = 79244
